$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column K (header + data rows) onto the new column L
# so the new "vat" column matches the existing look (font/alignment/etc.)
$ws.Range("K1:K5").Copy()
$ws.Range("L1:L5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Add new column L header "vat" and the vat values for each product row
$ws.Range("L1").Value = "vat"
$ws.Range("L2").Value = 5000
$ws.Range("L3").Value = 2000
$ws.Range("L4").Value = 1000
$ws.Range("L5").Value = 3000

# Update the active selection to match the saved workbook state
$ws.Range("L3").Select()
